$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Authors")

# Set the new keynote title in B3 (Shin Yoo's talk title)
$ws.Range("B3").Value = "SBST in the age of AI Systems - Challenges Ahead"

# Update the active selection to B9, matching the saved workbook state
$ws.Activate()
$ws.Range("B9").Select()
